$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 8349116
$ws.Cells.Item(32, 9).Value = 11499
$ws.Cells.Item(32, 10).Value = 10016640
$ws.Cells.Item(32, 11).Value = 11499
$ws.Cells.Item(32, 12).Value = 10016640
$ws.Cells.Item(32, 13).Value = -11173
$ws.Cells.Item(32, 14).Value = -10017292

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 132.44444
$ws.Cells.Item(33, 9).Value = 132.44444
$ws.Cells.Item(33, 11).Value = 132.44444
$ws.Cells.Item(33, 13).Value = 96.55556000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 30000
$ws.Cells.Item(74, 9).Value = 30000
$ws.Cells.Item(74, 11).Value = 30000
$ws.Cells.Item(74, 13).Value = -29064

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 30000
$ws.Cells.Item(77, 9).Value = 30000
$ws.Cells.Item(77, 11).Value = 150000
$ws.Cells.Item(77, 13).Value = -145320

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 11112651
$ws.Cells.Item(92, 9).Value = 6251544.5
$ws.Cells.Item(92, 11).Value = 6251544.5
$ws.Cells.Item(92, 13).Value = -6250296.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1536
$ws.Cells.Item(112, 10).Value = 2811
$ws.Cells.Item(112, 12).Value = 8433
$ws.Cells.Item(112, 14).Value = -10649

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1391.6842
$ws.Cells.Item(138, 9).Value = 1261.2941
$ws.Cells.Item(138, 10).Value = 2500
$ws.Cells.Item(138, 11).Value = 3783.8823
$ws.Cells.Item(138, 12).Value = 7500
$ws.Cells.Item(138, 13).Value = 1356.1177
$ws.Cells.Item(138, 14).Value = -17780

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1416.6666
$ws.Cells.Item(2, 9).Value = 1166.6666
$ws.Cells.Item(2, 11).Value = 1166.6666
$ws.Cells.Item(2, 13).Value = -1053.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9200.790999999999
$ws.Cells.Item(32, 9).Value = 3448.52
$ws.Cells.Item(32, 11).Value = 3448.52
$ws.Cells.Item(32, 13).Value = -3161.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 6759.8
$ws.Cells.Item(36, 9).Value = 2199.75
$ws.Cells.Item(36, 11).Value = 2199.75
$ws.Cells.Item(36, 13).Value = -1853.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 113937
$ws.Cells.Item(61, 9).Value = 3179.25
$ws.Cells.Item(61, 11).Value = 3179.25
$ws.Cells.Item(61, 13).Value = -2967.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 1658.3334
$ws.Cells.Item(63, 9).Value = 987.5
$ws.Cells.Item(63, 11).Value = 987.5
$ws.Cells.Item(63, 13).Value = -301.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 1658.3334
$ws.Cells.Item(66, 9).Value = 987.5
$ws.Cells.Item(66, 11).Value = 4937.5
$ws.Cells.Item(66, 13).Value = -1505.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 36889.758
$ws.Cells.Item(74, 9).Value = 54143.156
$ws.Cells.Item(74, 11).Value = 54143.156
$ws.Cells.Item(74, 13).Value = -53269.156

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 36889.758
$ws.Cells.Item(77, 9).Value = 54143.156
$ws.Cells.Item(77, 11).Value = 270715.78
$ws.Cells.Item(77, 13).Value = -266347.78

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1416.6666
$ws.Cells.Item(116, 9).Value = 1166.6666
$ws.Cells.Item(116, 11).Value = 1166.6666
$ws.Cells.Item(116, 13).Value = 1127.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 113937
$ws.Cells.Item(136, 9).Value = 3179.25
$ws.Cells.Item(136, 11).Value = 9537.75
$ws.Cells.Item(136, 13).Value = -6987.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1416.6666
$ws.Cells.Item(3, 9).Value = 1166.6666
$ws.Cells.Item(3, 11).Value = 1166.6666
$ws.Cells.Item(3, 13).Value = -1052.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 110395.93
$ws.Cells.Item(20, 9).Value = 131803.16
$ws.Cells.Item(20, 10).Value = 3359.8
$ws.Cells.Item(20, 11).Value = 131803.16
$ws.Cells.Item(20, 12).Value = 3359.8
$ws.Cells.Item(20, 13).Value = -131556.16
$ws.Cells.Item(20, 14).Value = -3853.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 83496
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 13).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(55, 8).Value = 17000
$ws.Cells.Item(55, 10).Value = 17000
$ws.Cells.Item(55, 12).Value = 17000
$ws.Cells.Item(55, 14).Value = -17630

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2752.9033
$ws.Cells.Item(122, 9).Value = 3069.3
$ws.Cells.Item(122, 10).Value = 2602.238
$ws.Cells.Item(122, 11).Value = 9207.900000000001
$ws.Cells.Item(122, 12).Value = 7806.714
$ws.Cells.Item(122, 13).Value = -6757.900000000001
$ws.Cells.Item(122, 14).Value = -12706.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 4529570
$ws.Cells.Item(134, 9).Value = 4764875
$ws.Cells.Item(134, 10).Value = 999999
$ws.Cells.Item(134, 11).Value = 14294625
$ws.Cells.Item(134, 12).Value = 2999997
$ws.Cells.Item(134, 13).Value = -14292090
$ws.Cells.Item(134, 14).Value = -3005067

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(138, 8).Value = 57782.855
$ws.Cells.Item(138, 9).Value = 50000
$ws.Cells.Item(138, 10).Value = 59080
$ws.Cells.Item(138, 11).Value = 50000
$ws.Cells.Item(138, 12).Value = 59080
$ws.Cells.Item(138, 13).Value = -44860
$ws.Cells.Item(138, 14).Value = -69360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 3383.5
$ws.Cells.Item(3, 9).Value = 3383.5
$ws.Cells.Item(3, 11).Value = 10150.5
$ws.Cells.Item(3, 13).Value = -10038.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 78307.30499999999
$ws.Cells.Item(14, 9).Value = 78307.30499999999
$ws.Cells.Item(14, 11).Value = 234921.915
$ws.Cells.Item(14, 13).Value = -234748.915

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 1680.4736
$ws.Cells.Item(121, 10).Value = 2465
$ws.Cells.Item(121, 12).Value = 7395
$ws.Cells.Item(121, 14).Value = -10015

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 1301.6666
$ws.Cells.Item(13, 10).Value = 1301.6666
$ws.Cells.Item(13, 12).Value = 1301.6666
$ws.Cells.Item(13, 14).Value = -1579.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 10337.167
$ws.Cells.Item(19, 9).Value = 11510.75
$ws.Cells.Item(19, 10).Value = 7990
$ws.Cells.Item(19, 11).Value = 11510.75
$ws.Cells.Item(19, 12).Value = 7990
$ws.Cells.Item(19, 13).Value = -11222.75
$ws.Cells.Item(19, 14).Value = -8566

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1958
$ws.Cells.Item(97, 9).Value = 1749.3334
$ws.Cells.Item(97, 10).Value = 2166.6667
$ws.Cells.Item(97, 11).Value = 1749.3334
$ws.Cells.Item(97, 12).Value = 2166.6667
$ws.Cells.Item(97, 13).Value = -1253.3334
$ws.Cells.Item(97, 14).Value = -3158.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1256.1818
$ws.Cells.Item(102, 9).Value = 1231.8
$ws.Cells.Item(102, 10).Value = 1500
$ws.Cells.Item(102, 11).Value = 1231.8
$ws.Cells.Item(102, 12).Value = 1500
$ws.Cells.Item(102, 13).Value = 390.2
$ws.Cells.Item(102, 14).Value = -4744

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 559.5333000000001
$ws.Cells.Item(107, 9).Value = 542.13635
$ws.Cells.Item(107, 10).Value = 607.375
$ws.Cells.Item(107, 11).Value = 542.13635
$ws.Cells.Item(107, 12).Value = 607.375
$ws.Cells.Item(107, 13).Value = 1377.86365
$ws.Cells.Item(107, 14).Value = -4447.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 11225399
$ws.Cells.Item(122, 9).Value = 18706504
$ws.Cells.Item(122, 10).Value = 3742.25
$ws.Cells.Item(122, 11).Value = 56119512
$ws.Cells.Item(122, 12).Value = 11226.75
$ws.Cells.Item(122, 13).Value = -56117062
$ws.Cells.Item(122, 14).Value = -16126.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3883.2
$ws.Cells.Item(132, 9).Value = 3051.4285
$ws.Cells.Item(132, 11).Value = 9154.2855
$ws.Cells.Item(132, 13).Value = -6624.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1115.2
$ws.Cells.Item(22, 9).Value = 1262.9286
$ws.Cells.Item(22, 10).Value = 927.1818
$ws.Cells.Item(22, 11).Value = 1262.9286
$ws.Cells.Item(22, 12).Value = 927.1818
$ws.Cells.Item(22, 13).Value = -967.9286
$ws.Cells.Item(22, 14).Value = -1517.1818

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 1115.2
$ws.Cells.Item(27, 9).Value = 1262.9286
$ws.Cells.Item(27, 10).Value = 927.1818
$ws.Cells.Item(27, 11).Value = 1262.9286
$ws.Cells.Item(27, 12).Value = 927.1818
$ws.Cells.Item(27, 13).Value = -1155.9286
$ws.Cells.Item(27, 14).Value = -1141.1818

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(34, 8).Value = 1066.5
$ws.Cells.Item(34, 9).Value = 1066.5
$ws.Cells.Item(34, 11).Value = 1066.5
$ws.Cells.Item(34, 13).Value = -894.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 9154567
$ws.Cells.Item(40, 9).Value = 4262.25
$ws.Cells.Item(40, 10).Value = 17288172
$ws.Cells.Item(40, 11).Value = 4262.25
$ws.Cells.Item(40, 12).Value = 17288172
$ws.Cells.Item(40, 13).Value = -4126.25
$ws.Cells.Item(40, 14).Value = -17288444

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 192762
$ws.Cells.Item(68, 9).Value = 211898.2
$ws.Cells.Item(68, 10).Value = 1400
$ws.Cells.Item(68, 11).Value = 211898.2
$ws.Cells.Item(68, 12).Value = 1400
$ws.Cells.Item(68, 13).Value = -211149.2
$ws.Cells.Item(68, 14).Value = -2898

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 192762
$ws.Cells.Item(71, 9).Value = 211898.2
$ws.Cells.Item(71, 10).Value = 1400
$ws.Cells.Item(71, 11).Value = 1059491
$ws.Cells.Item(71, 12).Value = 7000
$ws.Cells.Item(71, 13).Value = -1055747
$ws.Cells.Item(71, 14).Value = -14488

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2228.8333
$ws.Cells.Item(132, 9).Value = 2228.8333
$ws.Cells.Item(132, 11).Value = 6686.499899999999
$ws.Cells.Item(132, 13).Value = -4156.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 6200
$ws.Cells.Item(14, 9).Value = 6200
$ws.Cells.Item(14, 11).Value = 6200
$ws.Cells.Item(14, 13).Value = -6032

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1732.32
$ws.Cells.Item(122, 9).Value = 1177.0588
$ws.Cells.Item(122, 11).Value = 3531.1764
$ws.Cells.Item(122, 13).Value = -1081.1764
